# "sites partial, path fixes"
#
# Insert 4 new rows above the old row 30 (shifting the existing
# "partial/killed" block of rows down from 30-35 to 34-39), then
# populate the first three of those new rows (29-31) with three new
# site entries. Rows 32-33 are intentionally left blank, mirroring the
# existing blank separator row (29) that previously sat above the old
# block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 30:35 down to 34:39, opening up rows 29:33.
$ws.Rows("30:33").Insert()

# Column A first, then column B, so new shared-string entries are
# appended in the same column-major order as the source edit.
$ws.Range("A29").Value = "Ageneral"
$ws.Range("A30").Value = "DevlpAnes"
$ws.Range("A31").Value = "RCA Training"

$ws.Range("B29").Value = "Keyword-oriented site divided into a few categories, last updated 11/14. "
$ws.Range("B30").Value = "Site dedicated to anesthesia education in limited-resource countries. Hosts numerous pdf handouts, seminars, guidelines and textbook."
$ws.Range("B31").Value = "Robust British anesthesia educational modules for medical students and junior trainees. "

# Leave the cursor where the author left it.
$ws.Range("B31").Select() | Out-Null
